$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "Globo"
$ws.Range("B59").Value = "RJ TV 1"
$ws.Range("C59").Value = "Obras"
$ws.Range("D59").Value = "2025-04-05T12:06"
$ws.Range("E59").Value = "Negativo"
$ws.Range("F59").Value = "Ruas alagadas em Campos. Moradores do Novo Jockey reclamam que esse é um problema antigo. Repórter *ao vivo* do Cais da Lapa, Centro de Campos. Não chove no momento do ao vivo. Tempo nublado. Ontem, onde mais choveu foi Dores de Macabu. Equipe percorreu alguns bairros. Rua Manoel Viana Ribeiro, no jóquei, teve pontos de alagamento. Depoimentos de moradores. poder público não toma providências. Risco de doenças, como leptospirose. `"Poder público esqueceu Novo Jóckey`". Defesa Civil continua em alerta. Em nota, prefeitura informou que a Defesa Civil está ciente do problema no Novo Jóckey. Moradores podem solicitar para fazer bombeamento das águas da chuva. "

$ws.Range("A60").Value = "Globo"
$ws.Range("B60").Value = "RJ TV 1"
$ws.Range("C60").Value = "Trânsito"
$ws.Range("D60").Value = "2025-04-05T12:30"
$ws.Range("E60").Value = "Neutro"
$ws.Range("F60").Value = "Rolezinhos em Campos. Motociclistas fazem manobras e circulam sem capacete em vários pontos da cidade. Repórter *ao vivo*. Entre as vias, a Avenida 28 de Março. Imagens do rolezinho da madrugada de quarta-feira com flagrantes de várias irregularidades, entre elas, trafegar sem capacete e muito barulho. Em nota, Guarda Municipal informou que realiza ações em parceria com a Polícia Militar e realiza ações preventivas. Se o condutor for abordado pela Guarda, é encaminhado à delegacia. "
